# Actualización desde MV -datos-
# Appends new daily "Tasas de captación" rows (04-08-2021 .. 02-09-2021)
# and fixes the E148/F148 values for the last pre-existing row (03-08-2021).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing row 148 (03-08-2021): only columns E and F changed ---
$ws.Range("E148").Value = 0.28
$ws.Range("F148").Value = 0.26

# --- New rows 149-170 ---
$dates = @(
    "04-08-2021", "05-08-2021", "06-08-2021", "09-08-2021", "10-08-2021",
    "11-08-2021", "12-08-2021", "13-08-2021", "16-08-2021", "17-08-2021",
    "18-08-2021", "19-08-2021", "20-08-2021", "23-08-2021", "24-08-2021",
    "25-08-2021", "26-08-2021", "27-08-2021", "30-08-2021", "31-08-2021",
    "01-09-2021", "02-09-2021"
)

$data = @(
    @(0.07000000000000001, 0.13, 0.02, 0.13, 0.28),
    @(0.07000000000000001, 0.1,  0.01, 0.19, 0.24),
    @(0.07000000000000001, 0.1,  0.01, 0.15, 0.32),
    @(0.07000000000000001, 0.1,  0.01, 0.22, 0.33),
    @(0.07000000000000001, 0.12, 0.01, 0.19, 0.19),
    @(0.07000000000000001, 0.13, 0.01, 0.22, 0.35),
    @(0.07000000000000001, 0.11, 0.01, 0.17, 0.36),
    @(0.07000000000000001, 0.11, 0,    0.22, 0.24),
    @(0.07000000000000001, 0.1,  0.01, 0.16, 0.2),
    @(0.08,                0.12, 0.01, 0.25, 0.12),
    @(0.08,                0.14, 0.01, 0.18, 0.26),
    @(0.07000000000000001, 0.1,  0.01, 0.16, 0.33),
    @(0.08,                0.09, 0,    0.12, 0.29),
    @(0.07000000000000001, 0.07000000000000001, 0.01, 0.19, 0.21),
    @(0.07000000000000001, 0.12, 0,    0.23, 0.28),
    @(0.07000000000000001, 0.13, 0.01, 0.18, 0.24),
    @(0.08,                0.1,  0.01, 0.26, 0.25),
    @(0.08,                0.13, 0,    0.22, 0.26),
    @(0.07000000000000001, 0.1,  0.01, 0.18, 0.37),
    @(0.09,                0.13, 0.01, 0.18, 0.22),
    @(0.11,                0.1,  0.01, 0.11, 0.14),
    @(0.12,                0.13, 0.01, 0.22, 0.22)
)

$startRow = 149

# Column A holds dd-mm-yyyy text labels (e.g. "04-08-2021"). Labels whose
# day-of-month is <= 12 are ambiguous under default M-D-Y parsing and would
# otherwise be silently reinterpreted as serial dates, so force just those
# cells to Text before writing, so every label round-trips as a literal
# string (matching the rest of column A, which stores plain text dates).
$ws.Range("A149:A155").NumberFormat = "@"
$ws.Range("A169:A170").NumberFormat = "@"

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
    $ws.Cells.Item($row, 5).Value = $data[$i][3]
    $ws.Cells.Item($row, 6).Value = $data[$i][4]
}
